# Update the "Pais" worksheet with refreshed COVID-19 country/provincia data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 02:51"

# Row 4
$ws.Range("B4").Value = 8387798
$ws.Range("C4").Value = 44940
$ws.Range("D4").Value = 5457681
$ws.Range("E4").Value = 2705387
$ws.Range("G4").Value = 448
$ws.Range("H4").Value = 224730

# Row 5
$ws.Range("B5").Value = 7548238
$ws.Range("C5").Value = 55511
$ws.Range("D5").Value = 6659895
$ws.Range("E5").Value = 773701
$ws.Range("G5").Value = 578
$ws.Range("H5").Value = 114642

# Row 12
$ws.Range("B12").Value = 868675
$ws.Range("C12").Value = 3126
$ws.Range("D12").Value = 779779
$ws.Range("E12").Value = 55137
$ws.Range("G12").Value = 57
$ws.Range("H12").Value = 33759

# Row 21
$ws.Range("D21").Value = 291900
$ws.Range("E21").Value = 65215

# Row 133
$ws.Range("B133").Value = 5130
$ws.Range("C133").Value = 7
$ws.Range("D133").Value = 4944
$ws.Range("E133").Value = 77

# Row 194
$ws.Range("D194").Value = 121
$ws.Range("E194").Value = 26
